# Add a new slide at the end of the deck using the same "Title and Content" layout
# used throughout the rest of the deck.
$p = $ppt.ActivePresentation
$count = $p.Slides.Count
$s = $p.Slides.Add($count + 1, 2)

# ---- Title placeholder -------------------------------------------------
$titleTr = $s.Shapes.Item(1).TextFrame.TextRange
$titleTr.Text = "Repo "
$titleTr.Characters(1, 5).LanguageID = "en-US"

# ---- Body / content placeholder ----------------------------------------
$body = $s.Shapes.Item(2).TextFrame
$tr = $body.TextRange

$fullText = "Bit.ly/DriveTugas_CPP2023" + [char]13 + "https://bit.ly/BinusCPP2023 " + [char]13 + "Branch : Materi" + [char]13 + "Branch : Tugas " + [char]13 + "Branch : Sample" + [char]13
$tr.Text = $fullText

# Set the whole range to the Indonesian-tagged language used by the author,
# then restore the outline levels/hyperlinks/colors per run below.
$tr.LanguageID = "en-ID"

# Paragraph 1: "Bit.ly/DriveTugas_CPP2023" (hyperlink)
$run = $tr.Characters(1, 26)
$run.ActionSettings(1).Hyperlink.Address = "https://bit.ly/BinusCPP2023"
$run.Font.Color.RGB = 12673797

# Paragraph 2: "https://bit.ly/BinusCPP2023" (hyperlink) + " "
$run = $tr.Characters(28, 27)
$run.ActionSettings(1).Hyperlink.Address = "https://bit.ly/BinusCPP2023"
$run.Font.Color.RGB = 12673797

# Paragraph 3 (lvl 2 / IndentLevel 2): "Branch : " + "Materi"
$tr.Paragraphs(3).IndentLevel = 2

# Paragraph 4 (lvl 2 / IndentLevel 2): "Branch : " + "Tugas" + " "
$tr.Paragraphs(4).IndentLevel = 2

# Paragraph 5 (lvl 2 / IndentLevel 2): "Branch : Sample"
$tr.Paragraphs(5).IndentLevel = 2

Write-Output "Slide added: $($p.Slides.Count) total slides"
